$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the shared "short-url" text used in column B for every data row
#    (row 1 is the header "short-url", rows 2-375 hold the value "yAZD18")
$ws.Range("B2:B375").Value = "N8Wi39"

# 2. Apply the numeric corrections to the existing rows (no row shift yet)
$ws.Cells.Item(363, 14).Value = 5      # N363: 0 -> 5
$ws.Cells.Item(363, 20).Value = 10     # T363: 11 -> 10

$ws.Cells.Item(364, 14).Value = 361    # N364: 358 -> 361
$ws.Cells.Item(364, 15).Value = 55     # O364: 51 -> 55

$ws.Cells.Item(365, 14).Value = 30     # N365: 26 -> 30
$ws.Cells.Item(365, 15).Value = 5      # O365: 7 -> 5

$ws.Cells.Item(367, 14).Value = 4655   # N367: 4549 -> 4655
$ws.Cells.Item(367, 15).Value = 1208   # O367: 1071 -> 1208

$ws.Cells.Item(369, 20).Value = 38     # T369: 34 -> 38

$ws.Cells.Item(371, 14).Value = 127    # N371: 133 -> 127
$ws.Cells.Item(371, 15).Value = 9      # O371: 7 -> 9

# 3. Insert a brand-new row at position 373 (Sudan/Uganda/Zimbabwe shift down by one)
$ws.Rows.Item(373).Insert()

# 4. Populate the new "Stateless" row 373
$ws.Cells.Item(373, 1).Value = "1"
$ws.Cells.Item(373, 2).Value = "N8Wi39"
$ws.Cells.Item(373, 3).Value = "1"
$ws.Cells.Item(373, 4).Value = "372"
$ws.Cells.Item(373, 5).Value = "2024"
$ws.Cells.Item(373, 6).Value = "216"
$ws.Cells.Item(373, 7).Value = "Stateless"
$ws.Cells.Item(373, 8).Value = "STA"
$ws.Cells.Item(373, 9).Value = "XXA"
$ws.Cells.Item(373, 10).Value = "136"
$ws.Cells.Item(373, 11).Value = "Namibia"
$ws.Cells.Item(373, 12).Value = "NAM"
$ws.Cells.Item(373, 13).Value = "NAM"
$ws.Cells.Item(373, 14).Value = "0"
$ws.Cells.Item(373, 15).Value = "0"
$ws.Cells.Item(373, 16).Value = "0"
$ws.Cells.Item(373, 17).Value = "0"
$ws.Cells.Item(373, 18).Value = "0"
$ws.Cells.Item(373, 19).Value = "14796"
$ws.Cells.Item(373, 20).Value = "0"
$ws.Cells.Item(373, 21).Value = "-"
$ws.Cells.Item(373, 22).Value = "0"

# 5. Fix the "items" sequence number (column D) for the rows that shifted down
$ws.Cells.Item(374, 4).Value = "373"   # Sudan, was row 373
$ws.Cells.Item(375, 4).Value = "374"   # Uganda, was row 374
$ws.Cells.Item(376, 4).Value = "375"   # Zimbabwe, was row 375

# 6. Apply the numeric corrections to the rows that shifted down
$ws.Cells.Item(375, 14).Value = 14     # N375 (Uganda): 13 -> 14
$ws.Cells.Item(375, 15).Value = 11     # O375 (Uganda): 8 -> 11

$ws.Cells.Item(376, 14).Value = 43     # N376 (Zimbabwe): 42 -> 43
$ws.Cells.Item(376, 15).Value = 19     # O376 (Zimbabwe): 20 -> 19
